# Update "want to go" counts (column F) for a handful of events in the
# "展览" and "全部类型" worksheets, matching the regenerated data output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (Exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F3").Value = 16455
$wsExpo.Range("F8").Value = 9287
$wsExpo.Range("F13").Value = 228
$wsExpo.Range("F26").Value = 539
$wsExpo.Range("F30").Value = 86
$wsExpo.Range("F32").Value = 68
$wsExpo.Range("F33").Value = 273
$wsExpo.Range("F37").Value = 5725

# Sheet "全部类型" (All types) - same events appear at slightly different
# row numbers because this sheet also contains performance-type rows.
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F3").Value = 16455
$wsAll.Range("F8").Value = 9287
$wsAll.Range("F13").Value = 228
$wsAll.Range("F26").Value = 539
$wsAll.Range("F32").Value = 86
$wsAll.Range("F34").Value = 68
$wsAll.Range("F35").Value = 273
$wsAll.Range("F39").Value = 5725
